$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the floor (background) value across the full data range B2:K21
$ws.Range("B2:K21").Value = -18.94701190104432

# Override the specific non-floor score values per the updated PSSM
# Row 2
$ws.Range("C2").Value = 1.969037379524281

# Row 3
$ws.Range("I3").Value = 1.111610112394389

# Row 4
$ws.Range("C4").Value = 1.960348978989679
$ws.Range("D4").Value = 1.66419755246591
$ws.Range("F4").Value = 3.370755373905995
$ws.Range("H4").Value = 1.237603186773045
$ws.Range("J4").Value = 0.3113263935560134

# Row 5
$ws.Range("C5").Value = 1.639872772572907
$ws.Range("G5").Value = 2.669824523395834

# Row 7
$ws.Range("B7").Value = 2.430411005302832

# Row 8
$ws.Range("E8").Value = 1.807736173779686

# Row 9
$ws.Range("B9").Value = 3.868834367251823

# Row 10
$ws.Range("I10").Value = 1.583761136926735

# Row 11
$ws.Range("E11").Value = 2.914605510494364
$ws.Range("G11").Value = 2.855525368387781

# Row 13
$ws.Range("E13").Value = 2.521912877476256
$ws.Range("J13").Value = 1.694571558038179
$ws.Range("K13").Value = 4.321925382950077

# Row 14
$ws.Range("D14").Value = 1.506995143001117

# Row 15
$ws.Range("D15").Value = 1.810701195253237

# Row 16
$ws.Range("J16").Value = 1.861909080082767

# Row 17
$ws.Range("C17").Value = 2.141381046229649
$ws.Range("D17").Value = 1.836060668056448
$ws.Range("H17").Value = 2.093382846763566
$ws.Range("I17").Value = 2.162577119321069
$ws.Range("J17").Value = 2.620757146637425

# Row 18
$ws.Range("H18").Value = 2.079820165709351
$ws.Range("I18").Value = 2.085481180803241
$ws.Range("J18").Value = 2.520221978723572

# Row 19
$ws.Range("D19").Value = 2.026557252013891
$ws.Range("H19").Value = 1.69802462487687
$ws.Range("I19").Value = 1.832474224425403

# Row 20
$ws.Range("C20").Value = 1.110843685054411
$ws.Range("D20").Value = 1.50367682358525
$ws.Range("F20").Value = 3.271384896570347
$ws.Range("H20").Value = 1.589710556104011
$ws.Range("I20").Value = 1.355480884284313

# Row 21
$ws.Range("C21").Value = 1.325040320390079
$ws.Range("E21").Value = 1.68499727423301
$ws.Range("G21").Value = 2.677796492558543
$ws.Range("H21").Value = 1.532283855165169
